# daily auto push: 2026-01-16 22:35 UTC
# New data point for 2026/01/17 at hour 5 (ranking 201) is inserted right
# after the existing 2026/01/17 hour-2 row (row 666), pushing every
# following row down by one. Everything from the old row 666 onward moves
# to row+1; the sheet's used range grows from D707 to D708.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at 666 — Excel shifts rows 666..707 down to
# 667..708 for us (and grows the dimension/used-range automatically).
$ws.Rows(666).Insert()

# Column A holds a date-shaped value ("2026/01/17") that must stay literal
# text (matching every other date cell in the sheet), not get auto-parsed
# into a real Excel date serial. Entering it with a leading apostrophe
# forces text, then ClearFormats() drops the transient "quote prefix"
# style so the cell ends up plain/unstyled like its neighbours.
$ws.Cells(666, 1).Value = "'2026/01/17"
$ws.Cells(666, 1).ClearFormats()

$ws.Cells(666, 2).Value = "土"
$ws.Cells(666, 3).Value = 5
$ws.Cells(666, 4).Value = 201
